# Make the "Rural" electrolysis facilities (SUPH2*R1N / SUPH2*R2N) available
# in all regions on the "Availability" sheet, mirroring the existing
# "Central" (C1N/C2N) technology rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Availability")

# Each new technology gets its own row appended right after the existing
# data block (which ended at row 15), with every region column (F:J) set
# to 1 ("available") and the technology code written into column K.
$newTechs = @(
    @{ Row = 16; Code = "SUPH2ALKR1N" },
    @{ Row = 17; Code = "SUPH2ALKR2N" },
    @{ Row = 18; Code = "SUPH2PEMR1N" },
    @{ Row = 19; Code = "SUPH2PEMR2N" },
    @{ Row = 20; Code = "SUPH2SOER2N" }
)

foreach ($tech in $newTechs) {
    $r = $tech.Row

    # Region availability flags (DKISLBH, DKISL1, DKISL2, DKISL3, MAR)
    $ws.Range("F" + $r + ":J" + $r).Value = 1

    # Technology / Pset_PN code, using the same cell style ("Normal 42",
    # style index 24 in this workbook) as the other technology-code cells
    # in column K.
    $kCell = $ws.Range("K" + $r)
    $kCell.Style = "Normal 42"
    $kCell.Value = $tech.Code
}

# Restore the sheet's selection/zoom as left by the editor.
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 101
$ws.Range("L24").Select() | Out-Null
